$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking price strings
# (e.g. "63.42", "18.00") are not auto-converted to numbers by Excel,
# matching the source workbook where Price is stored as text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.130.29'
$ws.Range("D3").Value = '1.834.11'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '239.83'
$ws.Range("E5").Value = '  -1.96%  '
$ws.Range("D6").Value = '0.6636'
$ws.Range("E6").Value = '  -4.64%  '
$ws.Range("E8").Value = '  -3.68%  '
$ws.Range("D9").Value = '0.07346'
$ws.Range("E9").Value = '  -4.37%  '
$ws.Range("D10").Value = '22.72'
$ws.Range("E10").Value = '  -3.52%  '
$ws.Range("D11").Value = '0.07676'
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").Value = '1.841.35'
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").Value = '5.019'
$ws.Range("E13").Value = '  -2.56%  '
$ws.Range("D14").Value = '0.6745'
$ws.Range("E14").Value = '  -2.49%  '
$ws.Range("D15").Value = '86.39'
$ws.Range("E15").Value = '  -5.04%  '
$ws.Range("D16").Value = '6.111'
$ws.Range("E16").Value = '  -2.47%  '
$ws.Range("D17").Value = '29.136.57'
$ws.Range("E17").Value = '  -1.08%  '
$ws.Range("D18").Value = '0.000008234'
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("D19").Value = '228.39'
$ws.Range("E19").Value = '  -4.11%  '
$ws.Range("D20").Value = '12.49'
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '7.292'
$ws.Range("D23").Value = '0.9994'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '160.40'
$ws.Range("E24").Value = '  +0.32%  '
$ws.Range("D25").Value = '0.1417'
$ws.Range("E25").Value = '  -5.11%  '
$ws.Range("D26").Value = '8.653'
$ws.Range("E26").Value = '  -2.61%  '
$ws.Range("D27").Value = '18.00'
$ws.Range("E27").Value = '  -1.25%  '
$ws.Range("E28").Value = '  -2.13%  '
$ws.Range("D29").Value = '4.235'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '4.099'
$ws.Range("E30").Value = '  -1.09%  '
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").Value = '0.05327'
$ws.Range("E32").Value = '  +4.56%  '
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '0.7461'
$ws.Range("E34").Value = '  -3.52%  '
$ws.Range("D35").Value = '1.129'
$ws.Range("E35").Value = '  -1.65%  '
$ws.Range("D36").Value = '2.677'
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").Value = '1.319.62'
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("E38").Value = '  -3.83%  '
$ws.Range("D39").Value = '2.713'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '0.9233'
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("D41").Value = '6.023'
$ws.Range("E41").Value = '  +4.33%  '
$ws.Range("D42").Value = '0.9986'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").Value = '103.25'
$ws.Range("E43").Value = '  -2.74%  '
$ws.Range("D44").Value = '1.983.60'
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").Value = '0.5167'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000121'
$ws.Range("E46").Value = '  -3.18%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.758'
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '63.42'
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("B49").Value = 'XinFinNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D49").Value = '0.07611'
$ws.Range("E49").Value = '  +13.47%  '
$ws.Range("D50").Value = '9.262'
$ws.Range("E50").Value = '  -5.98%  '
$ws.Range("E51").Value = '  -0.11%  '
